# Populate the header row (S.No. / Name / Gender) that was added to the
# previously-empty worksheet, and leave the selection on C2 to match the
# saved view state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "S.No."
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Gender"

[void]$ws.Range("C2").Select()
